$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target values for columns D, L, M, N, O, P, S for rows 2..12
# (columns A,B,C,E,F,G,H,I,J,K,Q,R,T are identical across all rows and
# are left untouched)

$ws.Range("D2").Value = 44466
$ws.Range("M2").Value = 60

$ws.Range("D3").Value = 44503
$ws.Range("N3").Value = 30000
$ws.Range("O3").Value = 30000
$ws.Range("P3").Value = 30000
$ws.Range("S3").Value = 3000

$ws.Range("D4").Value = 44503
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 25000
$ws.Range("O4").Value = 25000
$ws.Range("P4").Value = 25000
$ws.Range("S4").Value = 2500

$ws.Range("D5").Value = 44435
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 40
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 20000
$ws.Range("S5").Value = 2000

$ws.Range("L6").Value = "Especial"
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 27000
$ws.Range("O6").Value = 27000
$ws.Range("P6").Value = 27000
$ws.Range("S6").Value = 2700

$ws.Range("D7").Value = 44517
$ws.Range("M7").Value = 30
$ws.Range("N7").Value = 25000
$ws.Range("O7").Value = 25000
$ws.Range("P7").Value = 25000
$ws.Range("S7").Value = 2500

$ws.Range("D8").Value = 44476
$ws.Range("M8").Value = 120
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 20000
$ws.Range("S8").Value = 2000

$ws.Range("D9").Value = 44434
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 20
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 20000
$ws.Range("S9").Value = 2000

$ws.Range("D11").Value = 44473
$ws.Range("M11").Value = 180

$ws.Range("D12").Value = 44511
$ws.Range("M12").Value = 120
$ws.Range("N12").Value = 28000
$ws.Range("O12").Value = 28000
$ws.Range("P12").Value = 28000
$ws.Range("S12").Value = 2800
